{"js": "// The activity guide text was retranslated: the old Slovak-leaning\n// sentence \"m\u00f4\u017eete pozorova\u0165 s\u00fahvezdie Ozvezdje Laboda 2022: ...\" is\n// replaced everywhere with the corrected Slovenian wording\n// \"2022: Datumi kampanje za opazovanje Ozvezdje Laboda: ...\".\n// This phrase appears identically in four separate paragraphs of the\n// document body, so find every occurrence and replace its text in place.\n\nconst oldText =\n  \"m\u00f4\u017eete pozorova\u0165 s\u00fahvezdie Ozvezdje Laboda 2022: 10.-19. avgust, 9.-18. september, 8.-17. oktober\";\nconst newText =\n  \"2022: Datumi kampanje za opazovanje Ozvezdje Laboda: 10.-19. avgust, 9.-18. september, 8.-17. oktober\";\n\nconst results = context.document.body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The activity guide text was retranslated: the old Slovak-leaning\n# sentence \"m\u00f4\u017eete pozorova\u0165 s\u00fahvezdie Ozvezdje Laboda 2022: ...\" is\n# replaced everywhere with the corrected Slovenian wording\n# \"2022: Datumi kampanje za opazovanje Ozvezdje Laboda: ...\".\n# This phrase appears identically in four separate paragraphs of the\n# document body, so run a single Find/Replace All over the whole story.\n\n$d = $word.ActiveDocument\n\n$oldText = \"m\u00f4\u017eete pozorova\u0165 s\u00fahvezdie Ozvezdje Laboda 2022: 10.-19. avgust, 9.-18. september, 8.-17. oktober\"\n$newText = \"2022: Datumi kampanje za opazovanje Ozvezdje Laboda: 10.-19. avgust, 9.-18. september, 8.-17. oktober\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n$find.Forward = $true\n$find.Wrap = 1          # wdFindContinue\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n"}
